$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.04553166666666666
$ws.Range("H2").Value = 0.136595
$ws.Range("I2").Value = 0.003566186696759492
$ws.Range("J2").Value = 0.003566186696759492
$ws.Range("M2").Value = 568.5612486666666
$ws.Range("N2").Value = 1705.683746
$ws.Range("O2").Value = 0.6737621253161296
$ws.Range("P2").Value = 0.6737621253161296
$ws.Range("Q2").Value = 25.88754125387444
$ws.Range("R2").Value = 232.98787128487
$ws.Range("S2").Value = 0.002402761528082783
$ws.Range("T2").Value = 0.002402761528082783

# Row 3
$ws.Range("G3").Value = 0.04553166666666666
$ws.Range("H3").Value = 0.136595
$ws.Range("I3").Value = 0.003566186696759492
$ws.Range("J3").Value = 0.003566186696759492
$ws.Range("M3").Value = 88.00803400000001
$ws.Range("O3").Value = 0.1042921587987053
$ws.Range("P3").Value = 0.1042921587987053
$ws.Range("Q3").Value = 4.007152468076667
$ws.Range("R3").Value = 36.06437221269
$ws.Range("S3").Value = 0.0003719253092842714
$ws.Range("T3").Value = 0.0003719253092842714

# Row 4
$ws.Range("G4").Value = 0.04553166666666666
$ws.Range("H4").Value = 0.136595
$ws.Range("I4").Value = 0.003566186696759492
$ws.Range("J4").Value = 0.003566186696759492
$ws.Range("M4").Value = 187.2912243333334
$ws.Range("N4").Value = 561.8736730000001
$ws.Range("O4").Value = 0.2219457158851651
$ws.Range("P4").Value = 0.2219457158851651
$ws.Range("Q4").Value = 8.527681595937223
$ws.Range("R4").Value = 76.749134363435
$ws.Range("S4").Value = 0.0007914998593924376
$ws.Range("T4").Value = 0.0007914998593924375

# Row 5
$ws.Range("I5").Value = 0.9955368819077489
$ws.Range("J5").Value = 0.995536881907749
$ws.Range("M5").Value = 568.5612486666666
$ws.Range("N5").Value = 1705.683746
$ws.Range("O5").Value = 0.6737621253161296
$ws.Range("P5").Value = 0.6737621253161296
$ws.Range("Q5").Value = 7226.767494690836
$ws.Range("R5").Value = 65040.90745221752
$ws.Range("S5").Value = 0.6707550453847576
$ws.Range("T5").Value = 0.6707550453847577

# Row 6
$ws.Range("I6").Value = 0.9955368819077489
$ws.Range("J6").Value = 0.995536881907749
$ws.Range("M6").Value = 88.00803400000001
$ws.Range("O6").Value = 0.1042921587987053
$ws.Range("P6").Value = 0.1042921587987053
$ws.Range("Q6").Value = 1118.636911809171
$ws.Range("R6").Value = 10067.73220628254
$ws.Range("S6").Value = 0.1038266905778909
$ws.Range("T6").Value = 0.1038266905778909

# Row 7
$ws.Range("I7").Value = 0.9955368819077489
$ws.Range("J7").Value = 0.995536881907749
$ws.Range("M7").Value = 187.2912243333334
$ws.Range("N7").Value = 561.8736730000001
$ws.Range("O7").Value = 0.2219457158851651
$ws.Range("P7").Value = 0.2219457158851651
$ws.Range("Q7").Value = 2380.588081279018
$ws.Range("R7").Value = 21425.29273151116
$ws.Range("S7").Value = 0.2209551459451004
$ws.Range("T7").Value = 0.2209551459451004

# Row 8
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.01145166666666667
$ws.Range("H8").Value = 0.034355
$ws.Range("I8").Value = 0.0008969313954915797
$ws.Range("J8").Value = 0.0008969313954915798
$ws.Range("M8").Value = 568.5612486666666
$ws.Range("N8").Value = 1705.683746
$ws.Range("O8").Value = 0.6737621253161296
$ws.Range("P8").Value = 0.6737621253161296
$ws.Range("Q8").Value = 6.510973899314443
$ws.Range("R8").Value = 58.59876509382999
$ws.Range("S8").Value = 0.0006043184032891687
$ws.Range("T8").Value = 0.0006043184032891687

# Row 9
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.01145166666666667
$ws.Range("H9").Value = 0.034355
$ws.Range("I9").Value = 0.0008969313954915797
$ws.Range("J9").Value = 0.0008969313954915798
$ws.Range("M9").Value = 88.00803400000001
$ws.Range("O9").Value = 0.1042921587987053
$ws.Range("P9").Value = 0.1042921587987053
$ws.Range("Q9").Value = 1.007838669356667
$ws.Range("R9").Value = 9.07054802421
$ws.Range("S9").Value = 0.00009354291153015221
$ws.Range("T9").Value = 0.00009354291153015223

# Row 10
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.01145166666666667
$ws.Range("H10").Value = 0.034355
$ws.Range("I10").Value = 0.0008969313954915797
$ws.Range("J10").Value = 0.0008969313954915798
$ws.Range("M10").Value = 187.2912243333334
$ws.Range("N10").Value = 561.8736730000001
$ws.Range("O10").Value = 0.2219457158851651
$ws.Range("P10").Value = 0.2219457158851651
$ws.Range("Q10").Value = 2.144796670657222
$ws.Range("R10").Value = 19.303170035915
$ws.Range("S10").Value = 0.0001990700806722588
$ws.Range("T10").Value = 0.0001990700806722588
